$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 132; this pushes existing rows 132-157 down to 133-158
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new record
$ws.Cells.Item(132, 1).Value = 11
$ws.Cells.Item(132, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(132, 3).Value = "Bíobío"
$ws.Cells.Item(132, 4).Value = 45275
$ws.Cells.Item(132, 5).Value = 8
$ws.Cells.Item(132, 6).Value = 100112037
$ws.Cells.Item(132, 7).Value = "Cebollín"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 80
$ws.Cells.Item(132, 11).Value = 5000
$ws.Cells.Item(132, 12).Value = 5000
$ws.Cells.Item(132, 13).Value = 5000
$ws.Cells.Item(132, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(132, 15).Value = "Región Metropolitana"
$ws.Cells.Item(132, 16).Value = 139
$ws.Cells.Item(132, 17).Value = 36
$ws.Cells.Item(132, 18).Value = "Hortaliza"
